$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "31.141.31"
$ws.Cells.Item(2, 5).Value = "  +4.47%  "

$ws.Cells.Item(3, 4).Value = "1.908.13"
$ws.Cells.Item(3, 5).Value = "  +1.81%  "

Set-TextValue ($ws.Cells.Item(4, 4)) "0.9947"
$ws.Cells.Item(4, 5).Value = "  -0.33%  "

Set-TextValue ($ws.Cells.Item(5, 4)) "247.04"
$ws.Cells.Item(5, 5).Value = "  +1.53%  "

Set-TextValue ($ws.Cells.Item(6, 4)) "0.9949"
$ws.Cells.Item(6, 5).Value = "  -0.31%  "

Set-TextValue ($ws.Cells.Item(7, 4)) "0.4953"
$ws.Cells.Item(7, 5).Value = "  +0.39%  "

Set-TextValue ($ws.Cells.Item(8, 4)) "0.2975"
$ws.Cells.Item(8, 5).Value = "  +2.57%  "

Set-TextValue ($ws.Cells.Item(9, 4)) "0.06764"
$ws.Cells.Item(9, 5).Value = "  +2.45%  "

$ws.Cells.Item(10, 4).Value = "1.886.74"
$ws.Cells.Item(10, 5).Value = "  +0.63%  "

Set-TextValue ($ws.Cells.Item(11, 4)) "16.91"
$ws.Cells.Item(11, 5).Value = "  +0.12%  "

Set-TextValue ($ws.Cells.Item(12, 4)) "0.07263"
$ws.Cells.Item(12, 5).Value = "  +1.43%  "

Set-TextValue ($ws.Cells.Item(13, 4)) "0.6814"
$ws.Cells.Item(13, 5).Value = "  +2.04%  "

Set-TextValue ($ws.Cells.Item(14, 4)) "5.070"
$ws.Cells.Item(14, 5).Value = "  +5.62%  "

Set-TextValue ($ws.Cells.Item(15, 4)) "88.35"
$ws.Cells.Item(15, 5).Value = "  +3.55%  "

$ws.Cells.Item(16, 4).Value = "30.981.56"
$ws.Cells.Item(16, 5).Value = "  +3.97%  "

Set-TextValue ($ws.Cells.Item(17, 4)) "0.000007967"
$ws.Cells.Item(17, 5).Value = "  +2.16%  "

Set-TextValue ($ws.Cells.Item(18, 4)) "0.9964"
$ws.Cells.Item(18, 5).Value = "  -0.17%  "

Set-TextValue ($ws.Cells.Item(19, 4)) "13.10"
$ws.Cells.Item(19, 5).Value = "  +3.02%  "

$ws.Cells.Item(20, 4).Value = "2.129.11"
$ws.Cells.Item(20, 5).Value = "  +0.48%  "

Set-TextValue ($ws.Cells.Item(21, 4)) "0.9981"
$ws.Cells.Item(21, 5).Value = "  +0.05%  "

Set-TextValue ($ws.Cells.Item(22, 4)) "4.842"
$ws.Cells.Item(22, 5).Value = "  +2.46%  "

Set-TextValue ($ws.Cells.Item(23, 4)) "5.992"
$ws.Cells.Item(23, 5).Value = "  +7.82%  "

Set-TextValue ($ws.Cells.Item(24, 4)) "163.32"
$ws.Cells.Item(24, 5).Value = "  +22.37%  "

Set-TextValue ($ws.Cells.Item(25, 4)) "9.311"
$ws.Cells.Item(25, 5).Value = "  +2.04%  "

Set-TextValue ($ws.Cells.Item(26, 4)) "152.28"
$ws.Cells.Item(26, 5).Value = "  +3.51%  "

Set-TextValue ($ws.Cells.Item(27, 4)) "17.50"
$ws.Cells.Item(27, 5).Value = "  +4.96%  "

Set-TextValue ($ws.Cells.Item(28, 4)) "1.941"
$ws.Cells.Item(28, 5).Value = "  +1.12%  "

Set-TextValue ($ws.Cells.Item(29, 4)) "1.429"
$ws.Cells.Item(29, 5).Value = "  +3.62%  "

Set-TextValue ($ws.Cells.Item(30, 4)) "4.286"
$ws.Cells.Item(30, 5).Value = "  +2.67%  "

Set-TextValue ($ws.Cells.Item(31, 4)) "0.08892"
$ws.Cells.Item(31, 5).Value = "  +2.77%  "

Set-TextValue ($ws.Cells.Item(32, 4)) "4.062"
$ws.Cells.Item(32, 5).Value = "  +4.09%  "

Set-TextValue ($ws.Cells.Item(33, 4)) "0.05267"
$ws.Cells.Item(33, 5).Value = "  +4.43%  "

Set-TextValue ($ws.Cells.Item(34, 4)) "0.7438"
$ws.Cells.Item(34, 5).Value = "  +5.37%  "

Set-TextValue ($ws.Cells.Item(35, 4)) "1.138"
$ws.Cells.Item(35, 5).Value = "  +2.86%  "

Set-TextValue ($ws.Cells.Item(36, 4)) "2.658"
$ws.Cells.Item(36, 5).Value = "  -0.30%  "

Set-TextValue ($ws.Cells.Item(37, 4)) "0.01859"
$ws.Cells.Item(37, 5).Value = "  +13.79%  "

Set-TextValue ($ws.Cells.Item(38, 4)) "2.737"
$ws.Cells.Item(38, 5).Value = "  +3.05%  "

Set-TextValue ($ws.Cells.Item(39, 4)) "2.222"
$ws.Cells.Item(39, 5).Value = "  +0.87%  "

Set-TextValue ($ws.Cells.Item(40, 4)) "0.9444"
$ws.Cells.Item(40, 5).Value = "  +1.51%  "

Set-TextValue ($ws.Cells.Item(43, 4)) "104.89"
$ws.Cells.Item(43, 5).Value = "  +2.48%  "

Set-TextValue ($ws.Cells.Item(44, 4)) "0.9983"
$ws.Cells.Item(44, 5).Value = "  +0.39%  "

Set-TextValue ($ws.Cells.Item(45, 4)) "7.654"
$ws.Cells.Item(45, 5).Value = "  +2.14%  "

Set-TextValue ($ws.Cells.Item(46, 4)) "0.1312"
$ws.Cells.Item(46, 5).Value = "  +4.67%  "

Set-TextValue ($ws.Cells.Item(47, 4)) "0.05812"
$ws.Cells.Item(47, 5).Value = "  +2.14%  "

Set-TextValue ($ws.Cells.Item(48, 4)) "33.23"
$ws.Cells.Item(48, 5).Value = "  +2.24%  "

Set-TextValue ($ws.Cells.Item(49, 4)) "8.504"
$ws.Cells.Item(49, 5).Value = "  +3.44%  "

Set-TextValue ($ws.Cells.Item(50, 4)) "0.3861"
$ws.Cells.Item(50, 5).Value = "  +4.67%  "

Set-TextValue ($ws.Cells.Item(51, 4)) "1.368"
$ws.Cells.Item(51, 5).Value = "  +2.47%  "

# Row 41: FraxShare -> TheSandbox
$ws.Cells.Item(41, 2).Value = "TheSandbox"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue ($ws.Cells.Item(41, 4)) "0.4322"
$ws.Cells.Item(41, 5).Value = "  +3.74%  "

# Row 42: TheSandbox -> FraxShare
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue ($ws.Cells.Item(42, 4)) "5.876"
$ws.Cells.Item(42, 5).Value = "  -3.16%  "
